$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.655945897102356
$ws.Range("B1").Value = 2.02492618560791
$ws.Range("C1").Value = 2.212754964828491
$ws.Range("D1").Value = 2.513699054718018
$ws.Range("E1").Value = 3.33297061920166
